# ---------------------------------------------------------------------------
# Refactor unit kerja references (and related letterhead fields) to use
# configuration-style ${...} placeholders instead of hard-coded text.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-RunBefore($range) {
    # Force a run boundary immediately before $range's start by toggling a
    # character-formatting property on/off; leaves no residual formatting
    # because the value is restored to what it already was.
    $b = $range.Bold
    $range.Bold = 1
    $range.Bold = $b
}

# ---------------------------------------------------------------------------
# 1 & 4) Both "BPS Kabupaten Hulu Sungai Tengah" occurrences in the body
#        become "BPS ${kabupaten}" (and "Kepala BPS ${kabupaten}").
#        "Kabupaten" -> lowercase "kabupaten", wrapped with new ${ / } runs.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 2; $i++) {
    $rng = $d.Content
    $found = $rng.Find.Execute("Kabupaten Hulu Sungai Tengah", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $wordStart = $rng.Start
    $wordEnd = $wordStart + 9          # length of "Kabupaten"
    $tailStart = $wordEnd              # start of " Hulu Sungai Tengah"
    $tailEnd = $rng.End

    # lower-case "Kabupaten" -> "kabupaten"
    $wordRng = $d.Range($wordStart, $wordEnd)
    $wordRng.Text = "kabupaten"

    # " Hulu Sungai Tengah" -> "}"
    $tailRng = $d.Range($tailStart, $tailEnd)
    $tailRng.Text = "}"

    # insert new "${" run immediately before "kabupaten"
    $insertPoint = $d.Range($wordStart, $wordStart)
    $insertPoint.InsertBefore('${')
    $dollarRange = $d.Range($wordStart, $wordStart + 2)
    Split-RunBefore($dollarRange)
}

# ---------------------------------------------------------------------------
# 2) "${unit_kerja" + "2" + "}" (three runs) collapse into a single run
#    "${unit_kerja2}" with identical text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("unit_kerja2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "unit_kerjaX2"   # force a real text change so the runs merge
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("unit_kerjaX2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $rng2.Text = "unit_kerja2"   # restore original text, now in one run
    }
}

# ---------------------------------------------------------------------------
# 3) "Barabai, ${tanggal_rapat}" -> "${ibukota}, ${tanggal_rapat}"
#    "Barabai" -> "ibukota", wrapped with new ${ / } runs.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Barabai, ${tanggal_rapat}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wordStart = $rng.Start
    $wordEnd = $wordStart + 7    # length of "Barabai"

    $wordRng = $d.Range($wordStart, $wordEnd)
    $wordRng.Text = "ibukota"

    $insertPoint = $d.Range($wordStart, $wordStart)
    $insertPoint.InsertBefore('${')
    $dollarRange = $d.Range($wordStart, $wordStart + 2)
    Split-RunBefore($dollarRange)

    $afterWordStart = $wordStart + 2 + 7   # past "${" + "ibukota"
    $insertPoint2 = $d.Range($afterWordStart, $afterWordStart)
    $insertPoint2.InsertBefore('}')
    $braceRange = $d.Range($afterWordStart, $afterWordStart + 1)
    Split-RunBefore($braceRange)
}

Write-Output $d.Content.Text
